# "Code Check In 8 Oct"
# Test Suite sheet:
#  - Row 20 (AdminSearch) now has a Description ("Done") and its Runmode
#    flips from "Y" to "N".
#  - A new row 21 is filled in for "AdminDashboard" (Runmode "N") — it
#    existed as a blank placeholder row before.
#  - A brand new row 22 is appended for "ReplayEvent" (Runmode "Y").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: AdminSearch
$ws.Range("B20").Value = "Done"
$ws.Range("C20").Value = "N"

# Row 21: AdminDashboard (was a blank row)
$ws.Range("A21").Value = "AdminDashboard"
$ws.Range("C21").Value = "N"

# Row 22: ReplayEvent (new row) - copy formatting down from row 21 first
# so the new row matches the rest of the table's styling, then fill values.
$ws.Range("A21:C21").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)

$ws.Range("A22").Value = "ReplayEvent"
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = "Y"

$ws.Range("B15").Select()
